# Update "想去人数" (column F) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 1346
    4  = 152
    7  = 99
    10 = 131
    11 = 4550
    12 = 6813
    18 = 4125
    22 = 2703
    26 = 353
    32 = 1021
    34 = 135
    40 = 80
    41 = 641
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
